$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.885.97"
$ws.Range("E2").Value = "  +0.04%  "

$ws.Range("D3").Value = "1.888.94"
$ws.Range("E3").Value = "  -0.30%  "

$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").Value = "'0.7680"
$ws.Range("E5").Value = "  -1.82%  "

$ws.Range("D6").Value = "'242.60"
$ws.Range("E6").Value = "  -0.60%  "

$ws.Range("D7").Value = "'1.001"
$ws.Range("E7").Value = "  +0.15%  "

$ws.Range("D8").Value = "'0.3133"
$ws.Range("E8").Value = "  -0.53%  "

$ws.Range("D9").Value = "'25.62"
$ws.Range("E9").Value = "  +0.61%  "

$ws.Range("D10").Value = "'0.07126"
$ws.Range("E10").Value = "  -2.74%  "

$ws.Range("D11").Value = "'0.08531"
$ws.Range("E11").Value = "  +4.91%  "

$ws.Range("D12").Value = "'0.7636"
$ws.Range("E12").Value = "  -0.61%  "

$ws.Range("D13").Value = "1.895.80"
$ws.Range("E13").Value = "  +2.57%  "

$ws.Range("D14").Value = "'5.366"
$ws.Range("E14").Value = "  -2.10%  "

$ws.Range("D15").Value = "'93.65"
$ws.Range("E15").Value = "  +0.28%  "

$ws.Range("E16").Value = "  -1.05%  "

$ws.Range("D17").Value = "29.937.41"
$ws.Range("E17").Value = "  +0.28%  "

$ws.Range("D18").Value = "'13.75"
$ws.Range("E18").Value = "  -1.54%  "

$ws.Range("D19").Value = "'244.40"
$ws.Range("E19").Value = "  -0.66%  "

$ws.Range("D20").Value = "'0.000007821"
$ws.Range("E20").Value = "  -0.67%  "

$ws.Range("D21").Value = "'0.9994"
$ws.Range("E21").Value = "  +0.04%  "

$ws.Range("D22").Value = "'8.012"
$ws.Range("E22").Value = "  -1.65%  "

$ws.Range("E23").Value = "  +0.08%  "

$ws.Range("D24").Value = "'0.1632"
$ws.Range("E24").Value = "  +3.03%  "

$ws.Range("D25").Value = "'9.381"
$ws.Range("E25").Value = "  -0.81%  "

$ws.Range("D26").Value = "'163.11"
$ws.Range("E26").Value = "  +0.48%  "

$ws.Range("E27").Value = "  -0.13%  "

$ws.Range("D28").Value = "'2.037"
$ws.Range("E28").Value = "  -0.01%  "

$ws.Range("E29").Value = "  +3.49%  "

$ws.Range("D30").Value = "'1.536"
$ws.Range("E30").Value = "  -0.52%  "

$ws.Range("D31").Value = "'4.506"
$ws.Range("E31").Value = "  +0.48%  "

$ws.Range("D32").Value = "'4.120"
$ws.Range("E32").Value = "  +0.86%  "

$ws.Range("D33").Value = "'0.05448"
$ws.Range("E33").Value = "  -2.78%  "

$ws.Range("D34").Value = "'1.243"
$ws.Range("E34").Value = "  -0.92%  "

$ws.Range("E36").Value = "  -0.14%  "

$ws.Range("D37").Value = "'2.705"
$ws.Range("E37").Value = "  +2.24%  "

$ws.Range("D38").Value = "'0.01950"
$ws.Range("E38").Value = "  +0.65%  "

$ws.Range("D39").Value = "'2.781"
$ws.Range("E39").Value = "  -0.16%  "

$ws.Range("D40").Value = "'0.4474"
$ws.Range("E40").Value = "  +0.23%  "

$ws.Range("D41").Value = "1.100.39"
$ws.Range("E41").Value = "  -3.68%  "

$ws.Range("D42").Value = "'73.17"
$ws.Range("E42").Value = "  -0.78%  "

$ws.Range("D43").Value = "'6.080"
$ws.Range("E43").Value = "  +1.77%  "

$ws.Range("D44").Value = "'0.8568"
$ws.Range("E44").Value = "  -0.18%  "

$ws.Range("E45").Value = "  +0.09%  "

$ws.Range("D46").Value = "'103.08"
$ws.Range("E46").Value = "  +0.99%  "

$ws.Range("D47").Value = "'7.680"
$ws.Range("E47").Value = "  +2.02%  "

$ws.Range("D48").Value = "'1.870"
$ws.Range("E48").Value = "  -1.60%  "

$ws.Range("E49").Value = "  -2.70%  "

$ws.Range("D50").Value = "2.035.15"
$ws.Range("E50").Value = "  +0.66%  "

$ws.Range("E51").Value = "  +0.33%  "
